# Generate Report for Handoff
# Update the localization-status report: the "fc40310c-...md" file has moved
# from "Handed back: in sync with en-US" to "Ready for handoff", with a new
# handoff timestamp recorded on the Overview sheet and on each language
# sheet's "Latest Handoff Datetime" column.

$wb = $excel.ActiveWorkbook

$statusText  = "Ready for handoff"

# ---- Overview sheet --------------------------------------------------
# Columns: A File Name | B zh-cn | C de-de | D Latest Handoff Date
# Row 3 = fc40310c-d3e2-430e-be7c-916464ffd6d6.md
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = $statusText
$ov.Range("C3").Value = $statusText
$ov.Range("D3").Value = "2016-03-24 21:01:45"

# ---- zh-cn sheet -------------------------------------------------------
# Columns: A Source File Name | B File Extension | C Status |
#          D Latest Handoff File | E Latest Handoff Datetime | ...
# Row 3 = fc40310c-d3e2-430e-be7c-916464ffd6d6.md
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $statusText
$zh.Range("E3").Value = "2016-03-24 21:01:37"

# ---- de-de sheet ---------------------------------------------------------
# Same column layout as zh-cn.
# Row 3 = fc40310c-d3e2-430e-be7c-916464ffd6d6.md
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $statusText
$de.Range("E3").Value = "2016-03-24 21:01:45"
